# "New crime data collected" — weekly refresh of the CompStat 18th Precinct
# workbook: bump the report's Volume/Number and week-covering dates by one
# week, and refresh the crime-count figures (and their derived % changes)
# for rows 15-27 and the Hate Crimes row (30).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Header text: "Volume 30   Number  25" -> "...Number  26"
#    and          "Report Covering the Week  6/19/2023  Through  6/25/2023"
#               -> "...6/26/2023  Through  7/2/2023"
# These are rich-text cells built from several runs; edit just the runs
# that changed (via Characters) so the surrounding runs/formatting and
# the rest of the shared string are left untouched.
# ---------------------------------------------------------------------

$volRange = $ws.Range("A8")
$volText = $volRange.Text
$volIdx = $volText.IndexOf("25") + 1
$volRange.Characters($volIdx, 2).Text = "26"

$weekRange = $ws.Range("C9")

$weekText = $weekRange.Text
$d1 = "6/19/2023"
$d1Idx = $weekText.IndexOf($d1) + 1
$weekRange.Characters($d1Idx, $d1.Length).Text = "6/26/2023"

$weekText2 = $weekRange.Text
$d2 = "6/25/2023"
$d2Idx = $weekText2.IndexOf($d2) + 1
$weekRange.Characters($d2Idx, $d2.Length).Text = "7/2/2023"

# ---------------------------------------------------------------------
# 2. Crime-count tables (rows 15-27, 30): new weekly/28-day/YTD/2yr figures
# ---------------------------------------------------------------------

# Row 15
$ws.Range("G15").Value = 3
$ws.Range("J15").Value = 10
$ws.Range("K15").Value = -50

# Row 16
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -66.666666666666
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = -60
$ws.Range("I16").Value = 53
$ws.Range("J16").Value = 72
$ws.Range("K16").Value = -26.388888888888
$ws.Range("L16").Value = 12.765957446808
$ws.Range("M16").Value = 3.921568627450
$ws.Range("N16").Value = -91.945288753799

# Row 17
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 13
$ws.Range("G17").Value = 8
$ws.Range("H17").Value = 62.5
$ws.Range("I17").Value = 93
$ws.Range("J17").Value = 84
$ws.Range("K17").Value = 10.714285714285
$ws.Range("L17").Value = 29.166666666666
$ws.Range("M17").Value = 40.909090909090
$ws.Range("N17").Value = -62.8

# Row 18
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 20
$ws.Range("H18").Value = -50
$ws.Range("I18").Value = 70
$ws.Range("J18").Value = 122
$ws.Range("K18").Value = -42.622950819672
$ws.Range("L18").Value = -23.913043478260
$ws.Range("M18").Value = -23.076923076923
$ws.Range("N18").Value = -93.560257589696

# Row 19
$ws.Range("C19").Value = 34
$ws.Range("E19").Value = -5.555555555555
$ws.Range("F19").Value = 146
$ws.Range("G19").Value = 145
$ws.Range("H19").Value = 0.689655172413
$ws.Range("I19").Value = 914
$ws.Range("J19").Value = 900
$ws.Range("K19").Value = 1.555555555555
$ws.Range("L19").Value = 83.534136546184
$ws.Range("M19").Value = 17.783505154639
$ws.Range("N19").Value = -74.462140262643

# Row 20
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 4
$ws.Range("G20").Value = 16
$ws.Range("H20").Value = -75
$ws.Range("I20").Value = 35
$ws.Range("J20").Value = 57
$ws.Range("K20").Value = -38.596491228070
$ws.Range("L20").Value = 20.689655172413
$ws.Range("M20").Value = 59.090909090909
$ws.Range("N20").Value = -86

# Row 21
$ws.Range("C21").Value = 41
$ws.Range("D21").Value = 50
$ws.Range("E21").Value = -18
$ws.Range("F21").Value = 177
$ws.Range("G21").Value = 202
$ws.Range("H21").Value = -12.376237623762
$ws.Range("I21").Value = 1170
$ws.Range("J21").Value = 1246
$ws.Range("K21").Value = -6.099518459069
$ws.Range("L21").Value = 56.208277703604
$ws.Range("M21").Value = 14.931237721021
$ws.Range("N21").Value = -79.993160054719

# Row 22
$ws.Range("C22").Value = 5
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 400
$ws.Range("F22").Value = 8
$ws.Range("G22").Value = 6
$ws.Range("H22").Value = 33.333333333333
$ws.Range("I22").Value = 36
$ws.Range("J22").Value = 32
$ws.Range("K22").Value = 12.5
$ws.Range("L22").Value = 33.333333333333
$ws.Range("M22").Value = 38.461538461538

# Row 24
$ws.Range("C24").Value = 60
$ws.Range("D24").Value = 48
$ws.Range("E24").Value = 25
$ws.Range("G24").Value = 231
$ws.Range("H24").Value = 2.164502164502
$ws.Range("I24").Value = 1331
$ws.Range("J24").Value = 1161
$ws.Range("K24").Value = 14.642549526270
$ws.Range("L24").Value = 54.408352668213
$ws.Range("M24").Value = 47.560975609756

# Row 25
$ws.Range("C25").Value = 17
$ws.Range("D25").Value = 14
$ws.Range("E25").Value = 21.428571428571
$ws.Range("F25").Value = 69
$ws.Range("G25").Value = 45
$ws.Range("H25").Value = 53.333333333333
$ws.Range("I25").Value = 352
$ws.Range("J25").Value = 272
$ws.Range("K25").Value = 29.411764705882
$ws.Range("L25").Value = 76
$ws.Range("M25").Value = 49.787234042553

# Row 26
$ws.Range("G26").Value = 3
$ws.Range("J26").Value = 18
$ws.Range("K26").Value = -33.333333333333
$ws.Range("L26").Value = -20

# Row 27
$ws.Range("C27").Value = 4
$ws.Range("E27").Value = 100
$ws.Range("F27").Value = 9
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = 28.571428571428
$ws.Range("I27").Value = 46
$ws.Range("J27").Value = 48
$ws.Range("K27").Value = -4.166666666666
$ws.Range("L27").Value = 31.428571428571

# Row 30 (Hate Crimes): Transit column (C) had 1 complaint, now has none —
# copy the "no data" text formatting from the Housing column (D, which
# already reads "0") so C30 becomes the same styled "0" placeholder.
$ws.Range("D30").Copy($ws.Range("C30"))
$ws.Range("F30").Value = 1
$ws.Range("H30").Value = 0
